$wb = $excel.ActiveWorkbook

# --- Sheet1: move selection from C3 back to A1 (it will lose the active/tabSelected tab) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
[void]$ws1.Range("A1").Select()

# --- Sheet2: populate with a simple "Message / Value / Number" data table ---
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A1").Value = "Message"
$ws2.Range("B1").Value = "Value"
$ws2.Range("C1").Value = "Number"

$ws2.Range("A2").Value = "kdhfdjks"
$ws2.Range("B2").Value = "regression34"
$ws2.Range("C2").Value = 24

$ws2.Range("A3").Value = "mvncxmn"
$ws2.Range("B3").Value = "regression24"
$ws2.Range("C3").Value = 32

$ws2.Range("A4").Value = "kjgfkgjjkhkj"
$ws2.Range("B4").Value = "regression25"
$ws2.Range("C4").Value = 9

# Sheet2 becomes the active/selected sheet with A1 selected
[void]$ws2.Range("A1").Select()
$ws2.Activate()
